$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Industries" policy column (H) for rows 25 through 78 from 1 to 0
$ws.Range("H25:H78").Value = 0
